$d = $word.ActiveDocument

# Locate the paragraph that currently reads "Modelo: {{marca}} {{modelo}}"
# (the list item describing the returned equipment's model).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Modelo:*{{marca}}*{{modelo}}*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $full = $target.Range

    # Remove the leading "Modelo: " label - this drops the two runs that
    # carried it ("Modelo:" in bold, then a plain space), leaving only the
    # run that holds the placeholder text untouched (including its rsid).
    $paraText = $full.Text
    $prefixLen = $paraText.IndexOf("{{marca}}")
    if ($prefixLen -lt 0) {
        $prefixLen = "Modelo: ".Length
    }
    $prefixRange = $d.Range($full.Start, $full.Start + $prefixLen)
    $prefixRange.Delete()

    # Replace the remaining placeholder text with the new wording.
    $remaining = $target.Range
    $remaining.Text = "Equipamento Devolvido: {{detalhes_equipamento}}"

    # Make sure the paragraph (including its paragraph mark) is bold, as in
    # the rest of the bold "label" runs on this template.
    $target.Range.Font.Bold = 1
    $target.Range.Font.BoldBi = 1
}
